# Update ASV_rank values in column G for rows 9-12 from 42 to 41
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G9").Value = 41
$ws.Range("G10").Value = 41
$ws.Range("G11").Value = 41
$ws.Range("G12").Value = 41
